$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Shrub" treatment group to "CSS" throughout the Tukey posthoc table.
# Because the group labels are re-sorted alphabetically ("CSS" now sorts
# before "Grassland" whereas "Shrub" sorted after it), the whole pairwise
# comparison table is regenerated in the new alphabetical order.

$data = @(
    @("CSS x Ambient",       "CSS x Reduced",       -0.0744, 0.9,    -1.2134, 1.0645,  $false),
    @("CSS x Ambient",       "Grassland x Ambient",  4.1294, 0.001,   2.9904, 5.2683,  $true),
    @("CSS x Ambient",       "Grassland x Reduced",  2.8087, 0.001,   1.6697, 3.9477,  $true),
    @("CSS x Reduced",       "Grassland x Ambient",  4.2038, 0.001,   3.0649, 5.3428,  $true),
    @("CSS x Reduced",       "Grassland x Reduced",  2.8831, 0.001,   1.7442, 4.0221,  $true),
    @("Grassland x Ambient", "Grassland x Reduced", -1.3207, 0.0168, -2.4596, -0.1817, $true)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
